# Apply WRI input-data fixes to the EoCtUH (Efficiency of Conversion to Usable
# Heat) workbook: swap the "Efficiency" header text for a clearer label,
# reorder the backing shared-string/value pair on the EoCtUH sheet, wrap the
# header text, grow the header row to fit, and leave the cursor parked on the
# header cell (while keeping "About" the active tab).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("About")
$ws3 = $wb.Worksheets.Item("EoCtUH")

# B1 used to read "Efficiency" (plain) and A2 used to read "District heat".
# The new layout swaps which string goes where AND renames the efficiency
# label to be more descriptive.
$ws3.Range("B1").Value = "Efficiency (dimensionless)"
$ws3.Range("A2").Value = "District heat"

# The header cell now wraps its (longer) text and the row grows to fit it.
$ws3.Range("B1").WrapText = $true
$ws3.Rows.Item(1).RowHeight = 45

# Leave the selection sitting on B1 for the EoCtUH sheet, but restore "About"
# as the active/visible tab afterwards (it was the active tab before the
# edit and stays that way).
$ws1.Activate()
$ws3.Range("B1").Select()
$ws1.Activate()
